# Re-order the header columns on the single worksheet.
# Original header row (A2:F2):
#   BusinessKey | Code | Name | Organization_ID | OrganizationTypeBusinessKey | ParentOrganization_ID
# New header row (A2:F2):
#   Organization_ID | BusinessKey | OrganizationTypeBusinessKey | Code | Name | ParentOrganization_ID

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Organization_ID"
$ws.Range("B2").Value = "BusinessKey"
$ws.Range("C2").Value = "OrganizationTypeBusinessKey"
$ws.Range("D2").Value = "Code"
$ws.Range("E2").Value = "Name"
$ws.Range("F2").Value = "ParentOrganization_ID"
